$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended after the existing row 21 (regcntr_id 10002..10010,
# usr_id 110021..110029, machine_id 10021..10029), following the same
# lang_code/is_active/cr_by/cr_dtimes/eff_dtimes pattern as the rows above.
$newRows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r++
}

# Active selection ends up on F14 as shown in the saved sheetView.
$ws.Range("F14").Select()

# Page setup: explicit portrait orientation (print resolution settings are
# stamped by Excel from the active printer driver and aren't reachable
# through the scripted object model here).
$ws.PageSetup.Orientation = 1
